$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, shifting existing rows 21-75 down to 22-76.
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new weekly data point.
$ws.Cells.Item(21, 1).Value = 10
$ws.Cells.Item(21, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(21, 3).Value = 'La Araucanía'
$ws.Cells.Item(21, 4).Value = 44459
$ws.Cells.Item(21, 5).Value = 9
$ws.Cells.Item(21, 6).Value = 100112031
$ws.Cells.Item(21, 7).Value = 'Poroto verde'
$ws.Cells.Item(21, 8).Value = 'Sin especificar'
$ws.Cells.Item(21, 9).Value = 'Primera'
$ws.Cells.Item(21, 10).Value = 20
$ws.Cells.Item(21, 11).Value = 35000
$ws.Cells.Item(21, 12).Value = 35000
$ws.Cells.Item(21, 13).Value = 35000
$ws.Cells.Item(21, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(21, 15).Value = 'Perú'
$ws.Cells.Item(21, 16).Value = 1400
$ws.Cells.Item(21, 17).Value = 25
$ws.Cells.Item(21, 18).Value = 'Hortaliza'
